# This edit shuffles the data rows (2-86) of the single data sheet: for
# each destination row, the values of columns D (Fecha), K (Variedad),
# L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion),
# R (Origen), S (Precio $/Kg) and T (Kg / unidad) are replaced with the
# values that used to live in another row (a pure permutation of the
# existing rows - no new data, nothing removed). Columns A, B, C, E, F,
# G, H, I, J are identical on every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values are copied from the source
# row's original contents into the destination row).
$rowMap = @{ 2 = 75; 3 = 76; 4 = 58; 5 = 59; 6 = 15; 7 = 16; 8 = 70; 9 = 86; 10 = 24; 11 = 25; 12 = 54; 13 = 8; 14 = 19; 15 = 39; 16 = 22; 17 = 38; 18 = 32; 19 = 41; 20 = 42; 21 = 33; 22 = 30; 23 = 73; 24 = 29; 25 = 72; 26 = 50; 27 = 9; 28 = 10; 29 = 78; 30 = 40; 31 = 27; 32 = 67; 33 = 84; 34 = 43; 35 = 3; 36 = 12; 37 = 13; 38 = 17; 39 = 82; 40 = 66; 41 = 57; 42 = 11; 43 = 51; 44 = 83; 45 = 69; 46 = 18; 47 = 14; 48 = 49; 49 = 62; 50 = 63; 51 = 36; 52 = 37; 53 = 77; 54 = 74; 55 = 28; 56 = 2; 57 = 55; 58 = 56; 59 = 4; 60 = 5; 61 = 6; 62 = 71; 63 = 7; 64 = 53; 65 = 31; 66 = 52; 67 = 64; 68 = 26; 69 = 65; 70 = 80; 71 = 81; 72 = 60; 73 = 44; 74 = 45; 75 = 21; 76 = 79; 77 = 61; 78 = 85; 79 = 34; 80 = 35; 81 = 23; 82 = 20; 83 = 68; 84 = 46; 85 = 47; 86 = 48 }

# Columns touched by the shuffle.
$cols = 4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20

$firstRow = 2
$lastRow = 86

# Snapshot the original contents of every affected column/row before
# writing anything back, so that source values are not clobbered while
# the destination rows are being updated in place.
$snapshot = @{}
foreach ($col in $cols) {
    $colData = @{}
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $colData[$r] = $ws.Cells.Item($r, $col).Value2
    }
    $snapshot[$col] = $colData
}

# Write the shuffled values back into the sheet.
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $snapshot[$col][$srcRow]
    }
}
